$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift each date in F2:F7 forward by one day (new chapters added, schedule moved)
$ws.Range("F2").Value = 44483
$ws.Range("F3").Value = 44482
$ws.Range("F4").Value = 44481
$ws.Range("F5").Value = 44480
$ws.Range("F6").Value = 44479
$ws.Range("F7").Value = 44478
